{"js": "// Apply the three text edits described by the diff using the Word\n// JavaScript API (Office.js) search/replace facility.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"...attached in the .java file. As...\" -> \"...attached in the dot java file. As...\"\nawait replaceOnce(\n  \"the .java file\",\n  \"the dot java file\"\n);\n\n// 2) Insert clause explaining \"average-case\" before \"selection sort will always out perform...\"\nawait replaceOnce(\n  \"In the \\u201Caverage-case\\u201D, selection sort will always out perform bubble sort.\",\n  \"In the \\u201Caverage-case\\u201D, which is taking the best case and worst case and divide it by 2, selection sort will always out perform bubble sort.\"\n);\n\n// 3) Rework the closing sentence about the merge sort being the best algorithm.\nawait replaceOnce(\n  \"The best sorting algorithm I believe is the merge sort and in code it would be recursion. \",\n  \"The best sorting algorithm however, I believe is the merge sort, which is recursion. \"\n);\n", "ps1": "# Apply the three text edits described by the diff using the Word COM\n# object model (Find/Replace against the document's main Range).\n#\n# Word Find.Execute positional signature:\n# (FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#  MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# 1) \"...attached in the .java file. As...\" -> \"...attached in the dot java file. As...\"\nReplace-Once \"the .java file\" \"the dot java file\"\n\n# 2) Insert clause explaining \"average-case\" before \"selection sort will always out perform...\"\nReplace-Once \"In the \u201caverage-case\u201d, selection sort will always out perform bubble sort.\" \"In the \u201caverage-case\u201d, which is taking the best case and worst case and divide it by 2, selection sort will always out perform bubble sort.\"\n\n# 3) Rework the closing sentence about the merge sort being the best algorithm.\nReplace-Once \"The best sorting algorithm I believe is the merge sort and in code it would be recursion. \" \"The best sorting algorithm however, I believe is the merge sort, which is recursion. \"\n"}
